# "Generate Report for Handoff"
# Updates the localization status report: the zh-cn / de-de files moved from
# "In Translation" to "Ready for handoff", so the Status columns (and the
# Overview sheet's per-language status columns) + the corresponding
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps are
# refreshed, and the now-wider status text causes those columns to re-autofit.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 13:14:13"

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-13 13:14:06"

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-13 13:14:13"

# Re-autofit the columns whose text just got wider ("Ready for handoff" vs
# "In Translation"), same as Excel does automatically when regenerating
# the report.
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(3).AutoFit() | Out-Null
